$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.329.29'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '3.215.75'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.215.08'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.71'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.35%  '
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').Value = '3.742.33'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').Value = '66.410.27'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('D18').Value = '3.216.25'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '506.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('E24').Value = '  -3.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.135'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +49.21%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  -5.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.44'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '502.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').Value = '0.0₃0776'
$ws.Range('E39').Value = '  +12.84%  '
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('E42').Value = '  +4.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').Value = '2.923.20'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.117'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
